$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws2 = $wb.Worksheets.Item("演出")
$ws4 = $wb.Worksheets.Item("全部类型")

# ---- 展览 ----
$ws1.Range("F2").Value = 2995
$ws1.Range("F3").Value = 6428
$ws1.Range("F4").Value = 2546
$ws1.Range("F6").Value = 648
$ws1.Range("F7").Value = 83
$ws1.Range("F9").Value = 3142
$ws1.Range("F10").Value = 365
$ws1.Range("F12").Value = 7654
$ws1.Range("F13").Value = 374
$ws1.Range("F16").Value = 259
$ws1.Range("F18").Value = 17
$ws1.Range("F19").Value = 485
$ws1.Range("F20").Value = 9382
$ws1.Range("F21").Value = 22
$ws1.Range("F26").Value = 28
$ws1.Range("F28").Value = 125
$ws1.Range("F30").Value = 128
$ws1.Range("F31").Value = 73
$ws1.Range("F32").Value = 118
$ws1.Range("F33").Value = 2621
$ws1.Range("F38").Value = 794
$ws1.Range("F39").Value = 3957
$ws1.Range("F43").Value = 107
$ws1.Range("F44").Value = 251
$ws1.Range("F45").Value = 48
$ws1.Range("F46").Value = 15
$ws1.Range("F47").Value = 67
$ws1.Range("F48").Value = 42
$ws1.Range("F49").Value = 64

# ---- 演出 ----
$ws2.Range("F6").Value = 272
$ws2.Range("F8").Value = 157
$ws2.Range("F16").Value = 9
$ws2.Range("F21").Value = 7
$ws2.Range("F23").Value = 9

# ---- 全部类型 ----
$ws4.Range("F3").Value = 2995
$ws4.Range("F5").Value = 272
$ws4.Range("F6").Value = 6428
$ws4.Range("F7").Value = 2546
$ws4.Range("F8").Value = 157
$ws4.Range("F10").Value = 648
$ws4.Range("F11").Value = 83
$ws4.Range("F13").Value = 3142
$ws4.Range("F14").Value = 365
$ws4.Range("B16").Value = "2024-04-27"
$ws4.Range("C16").Value = "北京·今泉爱夏  巡演"
$ws4.Range("D16").Value = "建国门外郎家园10号61幢一层A3-06、二层A3-06 East live"
$ws4.Range("E16").Value = "2024.04.27 20:00-04.27 21:30"
$ws4.Range("F16").Value = 48
$ws4.Range("G16").Value = 328
$ws4.Range("H16").Value = "https://show.bilibili.com/platform/detail.html?id=81889"
$ws4.Range("I16").Value = "//i1.hdslb.com/bfs/openplatform/202402/YJENeaUi1708313389899.jpeg"
$ws4.Range("C17").Value = "北京·2024中国爬宠狂欢节"
$ws4.Range("D17").Value = "丽泽天地购物中心 丽泽天地购物中心"
$ws4.Range("E17").Value = "2024.05.01 10:00-05.05 21:00"
$ws4.Range("F17").Value = 41
$ws4.Range("G17").Value = 30
$ws4.Range("H17").Value = "https://show.bilibili.com/platform/detail.html?id=83424"
$ws4.Range("I17").Value = "//i0.hdslb.com/bfs/openplatform/202403/kkbhUHKG1711524729125.jpeg"
$ws4.Range("C18").Value = "北京·IDO动漫游戏嘉年华45th"
$ws4.Range("F18").Value = 7654
$ws4.Range("G18").Value = 5
$ws4.Range("H18").Value = "https://show.bilibili.com/platform/detail.html?id=80645"
$ws4.Range("I18").Value = "//i0.hdslb.com/bfs/openplatform/202403/BIvjhmZq1709792042233.jpeg"
$ws4.Range("C19").Value = "北京·IDO动漫游戏嘉年华45th同人创作大会"
$ws4.Range("E19").Value = "2024.05.01 09:30-05.03 17:00"
$ws4.Range("F19").Value = 374
$ws4.Range("G19").Value = 85
$ws4.Range("H19").Value = "https://show.bilibili.com/platform/detail.html?id=82011"
$ws4.Range("I19").Value = "//i0.hdslb.com/bfs/openplatform/202402/2Aw7PvCg1708656416512.png"
$ws4.Range("C20").Value = "北京·动画电影《钢管公主》专场活动"
$ws4.Range("E20").Value = "2024.05.01 10:00-05.01 14:30"
$ws4.Range("F20").Value = 69
$ws4.Range("G20").Value = 528
$ws4.Range("H20").Value = "https://show.bilibili.com/platform/detail.html?id=83863"
$ws4.Range("I20").Value = "//i1.hdslb.com/bfs/openplatform/202404/oLIpAQh21712485244287.jpeg"
$ws4.Range("C21").Value = "北京·卡淘嘉年华·第三届球星卡交流会"
$ws4.Range("D21").Value = "亦庄荣昌东街6号 北京亦创国际会展中心"
$ws4.Range("E21").Value = "2024.05.01 09:30-05.03 17:00"
$ws4.Range("F21").Value = 109
$ws4.Range("H21").Value = "https://show.bilibili.com/platform/detail.html?id=82072"
$ws4.Range("I21").Value = "//i0.hdslb.com/bfs/openplatform/202402/XOTabMFt1708929919204.jpeg"
$ws4.Range("C22").Value = "北京·国乙同好嘉年华7th"
$ws4.Range("D22").Value = "北京国家会议中心 北京国家会议中心"
$ws4.Range("E22").Value = "2024.05.01 09:00-05.04 17:00"
$ws4.Range("F22").Value = 259
$ws4.Range("H22").Value = "https://show.bilibili.com/platform/detail.html?id=82391"
$ws4.Range("I22").Value = "//i1.hdslb.com/bfs/openplatform/202403/BGYIf9qe1709696198696.jpeg"
$ws4.Range("C23").Value = "北京·排球少年同好嘉年华2nd"
$ws4.Range("D23").Value = "天辰东路7号 北京国家会议中心"
$ws4.Range("E23").Value = "2024.05.01 09:30-05.04 17:00"
$ws4.Range("F23").Value = 17
$ws4.Range("G23").Value = 85
$ws4.Range("H23").Value = "https://show.bilibili.com/platform/detail.html?id=84070"
$ws4.Range("I23").Value = "//i2.hdslb.com/bfs/openplatform/202404/gkqcZt9X1712809702025.jpeg"
$ws4.Range("F24").Value = 9382
$ws4.Range("F25").Value = 22
$ws4.Range("F28").Value = 28
$ws4.Range("F30").Value = 125
$ws4.Range("F32").Value = 128
$ws4.Range("F33").Value = 73
$ws4.Range("F34").Value = 118
$ws4.Range("F35").Value = 2621
$ws4.Range("F38").Value = 794
$ws4.Range("F40").Value = 3957
$ws4.Range("F44").Value = 107
$ws4.Range("F45").Value = 251
$ws4.Range("F46").Value = 48
$ws4.Range("F47").Value = 67
$ws4.Range("F48").Value = 42
$ws4.Range("F49").Value = 64
